# Daily attendance processing - 2026-01-07 14:38:09
# Reorders the "Recorded By" (column G) values so that any leading
# "System" token is moved to the end of the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $value = $cell.Value2

    if ($value -ne $null -and $value -is [string] -and $value.StartsWith("System, ")) {
        $parts = $value -split ", "
        if ($parts[0] -eq "System") {
            $rest = $parts[1..($parts.Length - 1)]
            $newValue = ($rest -join ", ") + ", System"
            $cell.Value2 = $newValue
        }
    }
}
